$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws_LP1912 = $wb.Worksheets.Item("LP1912")
$ws_LP1912.Range("A2").Value = "Última actualización: 13:56:11"
$ws_LP1912.Range("A3").Value = "Total filas: 198"
$ws_LP1912.Range("C40").Value = "11_ETCHEVERRY"
$ws_LP1912.Range("C41").Value = "15_ABASTO"
$ws_LP1912.Range("A128").Value = "11:47:17"
$ws_LP1912.Range("C128").Value = "23_HERNANDEZ"
$ws_LP1912.Range("D128").Value = 45
$ws_LP1912.Range("A129").Value = "10:37:52"
$ws_LP1912.Range("C129").Value = "14_ABASTO"
$ws_LP1912.Range("D129").Value = 115
$ws_LP1912.Range("C139").Value = "16_SANTA ANA"
$ws_LP1912.Range("C141").Value = "14_ABASTO"
$ws_LP1912.Range("A167").Value = "12:45:56"
$ws_LP1912.Range("C167").Value = "23_HERNANDEZ"
$ws_LP1912.Range("D167").Value = 76
$ws_LP1912.Range("A168").Value = "12:11:52"
$ws_LP1912.Range("C168").Value = "10_OLMOS"
$ws_LP1912.Range("D168").Value = 110
$ws_LP1912.Range("A178").Value = "13:56:11"
$ws_LP1912.Range("B178").Value = "14:18"
$ws_LP1912.Range("C178").Value = "11_ETCHEVERRY"
$ws_LP1912.Range("D178").Value = 22
$ws_LP1912.Range("A179").Value = "12:45:56"
$ws_LP1912.Range("B179").Value = "14:27"
$ws_LP1912.Range("D179").Value = 102
$ws_LP1912.Range("A180").Value = "13:41:54"
$ws_LP1912.Range("B180").Value = "14:28"
$ws_LP1912.Range("C180").Value = "16_SANTA ANA"
$ws_LP1912.Range("D180").Value = 47
$ws_LP1912.Range("A181").Value = "12:45:56"
$ws_LP1912.Range("B181").Value = "14:31"
$ws_LP1912.Range("D181").Value = 106
$ws_LP1912.Range("A182").Value = "12:33:21"
$ws_LP1912.Range("B182").Value = "14:32"
$ws_LP1912.Range("C182").Value = "14X44_ABASTO"
$ws_LP1912.Range("D182").Value = 119
$ws_LP1912.Range("A183").Value = "12:45:56"
$ws_LP1912.Range("B183").Value = "14:33"
$ws_LP1912.Range("D183").Value = 108
$ws_LP1912.Range("A184").Value = "13:14:29"
$ws_LP1912.Range("B184").Value = "14:34"
$ws_LP1912.Range("C184").Value = "215C_EL PATO"
$ws_LP1912.Range("D184").Value = 80
$ws_LP1912.Range("A185").Value = "12:45:56"
$ws_LP1912.Range("B185").Value = "14:39"
$ws_LP1912.Range("C185").Value = "16_P MOR-SANTA ANA"
$ws_LP1912.Range("D185").Value = 114
$ws_LP1912.Range("B186").Value = "14:47"
$ws_LP1912.Range("C186").Value = "215B_EL PATO"
$ws_LP1912.Range("D186").Value = 115
$ws_LP1912.Range("A187").Value = "13:56:11"
$ws_LP1912.Range("B187").Value = "14:48"
$ws_LP1912.Range("C187").Value = "215B_EL PATO"
$ws_LP1912.Range("D187").Value = 52
$ws_LP1912.Range("A188").Value = "12:52:52"
$ws_LP1912.Range("B188").Value = "14:51"
$ws_LP1912.Range("C188").Value = "16_SANTA ANA"
$ws_LP1912.Range("D188").Value = 119
$ws_LP1912.Range("A189").Value = "13:41:54"
$ws_LP1912.Range("B189").Value = "14:51"
$ws_LP1912.Range("C189").Value = "23_HERNANDEZ"
$ws_LP1912.Range("D189").Value = 70
$ws_LP1912.Range("B190").Value = "14:53"
$ws_LP1912.Range("C190").Value = "215_EL PELIGRO"
$ws_LP1912.Range("D190").Value = 72
$ws_LP1912.Range("B191").Value = "14:54"
$ws_LP1912.Range("C191").Value = "215_EL PELIGRO"
$ws_LP1912.Range("D191").Value = 100
$ws_LP1912.Range("A192").Value = "13:41:54"
$ws_LP1912.Range("B192").Value = "15:01"
$ws_LP1912.Range("C192").Value = "10_OLMOS"
$ws_LP1912.Range("D192").Value = 80
$ws_LP1912.Range("B193").Value = "15:02"
$ws_LP1912.Range("C193").Value = "10_OLMOS"
$ws_LP1912.Range("D193").Value = 108
$ws_LP1912.Range("A194").Value = "13:14:29"
$ws_LP1912.Range("B194").Value = "15:12"
$ws_LP1912.Range("D194").Value = 118
$ws_LP1912.Range("A195").Value = "13:14:29"
$ws_LP1912.Range("B195").Value = "15:13"
$ws_LP1912.Range("C195").Value = "17X38_ROMERO"
$ws_LP1912.Range("D195").Value = 119
$ws_LP1912.Range("B196").Value = "15:14"
$ws_LP1912.Range("C196").Value = "14_ABASTO"
$ws_LP1912.Range("D196").Value = 93
# New rows for LP1912
$ws_LP1912.Range("A197").Value = "13:56:11"
$ws_LP1912.Range("B197").Value = "15:18"
$ws_LP1912.Range("C197").Value = "14_ABASTO"
$ws_LP1912.Range("D197").Value = 82
$ws_LP1912.Range("E197").Value = "LP1912"
$ws_LP1912.Range("A198").Value = "13:41:54"
$ws_LP1912.Range("B198").Value = "15:32"
$ws_LP1912.Range("C198").Value = "11_ETCHEVERRY"
$ws_LP1912.Range("D198").Value = 111
$ws_LP1912.Range("E198").Value = "LP1912"
$ws_LP1912.Range("A199").Value = "13:41:54"
$ws_LP1912.Range("B199").Value = "15:33"
$ws_LP1912.Range("C199").Value = "215C_EL PATO"
$ws_LP1912.Range("D199").Value = 112
$ws_LP1912.Range("E199").Value = "LP1912"
$ws_LP1912.Range("A200").Value = "13:56:11"
$ws_LP1912.Range("B200").Value = "15:34"
$ws_LP1912.Range("C200").Value = "215C_EL PATO"
$ws_LP1912.Range("D200").Value = 98
$ws_LP1912.Range("E200").Value = "LP1912"
$ws_LP1912.Range("A201").Value = "13:56:11"
$ws_LP1912.Range("B201").Value = "15:42"
$ws_LP1912.Range("C201").Value = "11_ETCHEVERRY"
$ws_LP1912.Range("D201").Value = 106
$ws_LP1912.Range("E201").Value = "LP1912"
$ws_LP1912.Range("A202").Value = "13:56:11"
$ws_LP1912.Range("B202").Value = "15:53"
$ws_LP1912.Range("C202").Value = "15X38_ABASTO"
$ws_LP1912.Range("D202").Value = 117
$ws_LP1912.Range("E202").Value = "LP1912"
$ws_LP1912.Range("A203").Value = "13:56:11"
$ws_LP1912.Range("B203").Value = "15:53"
$ws_LP1912.Range("C203").Value = "16_P MOR-SANTA ANA"
$ws_LP1912.Range("D203").Value = 117
$ws_LP1912.Range("E203").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws_LP1912_215 = $wb.Worksheets.Item("LP1912-215")
$ws_LP1912_215.Range("A2").Value = "Última actualización: 13:56:11"
$ws_LP1912_215.Range("A3").Value = "Total filas: 34"
$ws_LP1912_215.Range("A35").Value = "13:56:11"
$ws_LP1912_215.Range("B35").Value = "14:48"
$ws_LP1912_215.Range("C35").Value = "215B_EL PATO"
$ws_LP1912_215.Range("D35").Value = 52
$ws_LP1912_215.Range("A36").Value = "13:41:54"
$ws_LP1912_215.Range("B36").Value = "14:53"
$ws_LP1912_215.Range("D36").Value = 72
$ws_LP1912_215.Range("A37").Value = "13:14:29"
$ws_LP1912_215.Range("B37").Value = "14:54"
$ws_LP1912_215.Range("C37").Value = "215_EL PELIGRO"
$ws_LP1912_215.Range("D37").Value = 100
# New rows for LP1912-215
$ws_LP1912_215.Range("A38").Value = "13:41:54"
$ws_LP1912_215.Range("B38").Value = "15:33"
$ws_LP1912_215.Range("C38").Value = "215C_EL PATO"
$ws_LP1912_215.Range("D38").Value = 112
$ws_LP1912_215.Range("E38").Value = "LP1912"
$ws_LP1912_215.Range("A39").Value = "13:56:11"
$ws_LP1912_215.Range("B39").Value = "15:34"
$ws_LP1912_215.Range("C39").Value = "215C_EL PATO"
$ws_LP1912_215.Range("D39").Value = 98
$ws_LP1912_215.Range("E39").Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws_6203_6173 = $wb.Worksheets.Item("6203-6173")
$ws_6203_6173.Range("A2").Value = "Última actualización: 13:56:11"
$ws_6203_6173.Range("A3").Value = "Total filas: 28"
$ws_6203_6173.Range("A19").Value = "08:52:50"
$ws_6203_6173.Range("C19").Value = "215A_LA PLATA"
$ws_6203_6173.Range("D19").Value = 98
$ws_6203_6173.Range("A20").Value = "08:37:25"
$ws_6203_6173.Range("C20").Value = "215B_LP-P MOR-1 Y 57"
$ws_6203_6173.Range("D20").Value = 113
$ws_6203_6173.Range("A29").Value = "13:56:11"
$ws_6203_6173.Range("B29").Value = "14:27"
$ws_6203_6173.Range("D29").Value = 31
$ws_6203_6173.Range("A30").Value = "13:41:54"
$ws_6203_6173.Range("B30").Value = "14:28"
$ws_6203_6173.Range("D30").Value = 47
$ws_6203_6173.Range("A31").Value = "13:14:29"
$ws_6203_6173.Range("B31").Value = "14:33"
$ws_6203_6173.Range("C31").Value = "215C_LA PLATA"
$ws_6203_6173.Range("D31").Value = 79
$ws_6203_6173.Range("E31").Value = "L6203"
# New rows for 6203-6173
$ws_6203_6173.Range("A32").Value = "13:41:54"
$ws_6203_6173.Range("B32").Value = "15:21"
$ws_6203_6173.Range("C32").Value = "215A_LA PLATA"
$ws_6203_6173.Range("D32").Value = 100
$ws_6203_6173.Range("E32").Value = "L6173"
$ws_6203_6173.Range("A33").Value = "13:56:11"
$ws_6203_6173.Range("B33").Value = "15:22"
$ws_6203_6173.Range("C33").Value = "215A_LA PLATA"
$ws_6203_6173.Range("D33").Value = 86
$ws_6203_6173.Range("E33").Value = "L6173"
